$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.057.68"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "3.389.36"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'573.40"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'137.38"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D8").Value = "3.387.29"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.470"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").Value = "'7.64"
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("D13").Value = "3.967.53"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "'26.69"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.386.94"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "61.120.17"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "'13.90"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "'375.85"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("B23").Value = "WrappedeETH"
$ws.Range("C23").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D23").Value = "3.514.45"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "'0.551"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "'70.86"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "'0.0000124"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").Value = "'1.62"
$ws.Range("E28").Value = "  -6.62%  "
$ws.Range("D29").Value = "'0.174"
$ws.Range("E29").Value = "  +8.50%  "
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").Value = "'7.40"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").Value = "'8.05"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "'5.12"
$ws.Range("E36").Value = "  -3.01%  "
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").Value = "'6.83"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'165.11"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Value = "'26.04"
$ws.Range("E41").Value = "  +5.88%  "
$ws.Range("D42").Value = "'1.76"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'0.774"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'41.91"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("D48").Value = "2.517.70"
$ws.Range("E48").Value = "  +7.39%  "
$ws.Range("D49").Value = "'23.57"
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("E51").Value = "  +3.14%  "
